$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3158
$ws.Range("I3").Value = 3246
$ws.Range("I4").Value = 766
$ws.Range("I6").Value = 3696
$ws.Range("I7").Value = 11165

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 102
$ws.Range("I5").Value = 37
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 377
$ws.Range("I8").Value = 704
$ws.Range("I11").Value = 176
$ws.Range("I14").Value = 58
$ws.Range("I16").Value = 30
$ws.Range("I19").Value = 298
$ws.Range("I20").Value = 285
$ws.Range("I22").Value = 28
$ws.Range("I23").Value = 101
$ws.Range("I25").Value = 51
$ws.Range("I27").Value = 98
$ws.Range("I29").Value = 729
$ws.Range("I30").Value = 39
$ws.Range("I31").Value = 100
$ws.Range("I33").Value = 504
$ws.Range("I36").Value = 152
$ws.Range("I37").Value = 362
$ws.Range("I40").Value = 19
$ws.Range("I42").Value = 394
$ws.Range("I48").Value = 130
$ws.Range("I49").Value = 85
$ws.Range("I51").Value = 101
$ws.Range("I52").Value = 239
$ws.Range("I53").Value = 122
$ws.Range("I54").Value = 249
$ws.Range("I55").Value = 122
$ws.Range("I63").Value = 39
$ws.Range("I65").Value = 247
$ws.Range("I67").Value = 448
$ws.Range("I72").Value = 39
$ws.Range("I74").Value = 26
$ws.Range("I76").Value = 172
$ws.Range("I79").Value = 284
$ws.Range("I83").Value = 228
$ws.Range("I85").Value = 514
$ws.Range("I88").Value = 103
$ws.Range("I89").Value = 123
$ws.Range("I90").Value = 138
$ws.Range("I91").Value = 133
$ws.Range("I96").Value = 128
$ws.Range("I98").Value = 69
$ws.Range("I99").Value = 204
$ws.Range("I101").Value = 11165

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 132
$ws.Range("I3").Value = 205
$ws.Range("I7").Value = 514

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 239

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 78
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 221
$ws.Range("I7").Value = 704

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 54
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 111
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 377

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 32
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 116
$ws.Range("I4").Value = 27
$ws.Range("I7").Value = 362

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 157
$ws.Range("I7").Value = 448

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 80
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 123
$ws.Range("I6").Value = 164
$ws.Range("I7").Value = 504

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I3").Value = 8
$ws.Range("I4").Value = 8
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I6").Value = 196
$ws.Range("I7").Value = 729

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 71
$ws.Range("I7").Value = 172

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 135
$ws.Range("I4").Value = 34
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 394

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 90
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 285

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I3").Value = 19
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 43
$ws.Range("I3").Value = 28
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("I3").Value = 8
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 26
